$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 367.85715
$ws.Range("I12").Value = 409.66666
$ws.Range("K12").Value = 409.66666
$ws.Range("M12").Value = -239.66666
$ws.Range("H33").Value = 2914.8
$ws.Range("J33").Value = 3412
$ws.Range("L33").Value = 3412
$ws.Range("N33").Value = -3870
$ws.Range("H54").Value = 22000
$ws.Range("I54").Value = 15000
$ws.Range("K54").Value = 15000
$ws.Range("M54").Value = -14514
$ws.Range("H80").Value = 2150.5
$ws.Range("I80").Value = 419
$ws.Range("K80").Value = 1257
$ws.Range("M80").Value = -259
$ws.Range("H83").Value = 2150.5
$ws.Range("I83").Value = 419
$ws.Range("K83").Value = 3771
$ws.Range("M83").Value = 1221
$ws.Range("H86").Value = 4556.8
$ws.Range("I86").Value = 4596
$ws.Range("J86").Value = 4400
$ws.Range("K86").Value = 4596
$ws.Range("L86").Value = 4400
$ws.Range("M86").Value = -3473
$ws.Range("N86").Value = -6646
$ws.Range("H89").Value = 4556.8
$ws.Range("I89").Value = 4596
$ws.Range("J89").Value = 4400
$ws.Range("K89").Value = 22980
$ws.Range("L89").Value = 22000
$ws.Range("M89").Value = -17364
$ws.Range("N89").Value = -33232
$ws.Range("H99").Value = 845
$ws.Range("I99").Value = 845
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2535
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1037
$ws.Range("N99").Value = $null
$ws.Range("H116").Value = 4696.5
$ws.Range("I116").Value = 4629
$ws.Range("K116").Value = 4629
$ws.Range("M116").Value = -1187
$ws.Range("H132").Value = 58829424
$ws.Range("I132").Value = 62506176
$ws.Range("J132").Value = 1399
$ws.Range("K132").Value = 187518528
$ws.Range("L132").Value = 4197
$ws.Range("M132").Value = -187515998
$ws.Range("N132").Value = -9257
$ws.Range("H136").Value = 226666.33
$ws.Range("J136").Value = 226666.33
$ws.Range("L136").Value = 226666.33
$ws.Range("N136").Value = -236866.33

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 7169.3335
$ws.Range("H24").Value = 17500
$ws.Range("J24").Value = 17500
$ws.Range("L24").Value = 17500
$ws.Range("N24").Value = -18248
$ws.Range("H38").Value = 6746.2
$ws.Range("I38").Value = 6746.2
$ws.Range("K38").Value = 6746.2
$ws.Range("M38").Value = -6279.2
$ws.Range("H45").Value = 1390.6666
$ws.Range("I45").Value = 1244.3636
$ws.Range("K45").Value = 1244.3636
$ws.Range("M45").Value = -867.3635999999999
$ws.Range("H74").Value = 2003.1765
$ws.Range("I74").Value = 1815.9375
$ws.Range("K74").Value = 1815.9375
$ws.Range("M74").Value = -941.9375
$ws.Range("H77").Value = 2003.1765
$ws.Range("I77").Value = 1815.9375
$ws.Range("K77").Value = 9079.6875
$ws.Range("M77").Value = -4711.6875
$ws.Range("H81").Value = 55998
$ws.Range("J81").Value = 55998
$ws.Range("L81").Value = 55998
$ws.Range("N81").Value = -57994
$ws.Range("H84").Value = 55998
$ws.Range("J84").Value = 55998
$ws.Range("L84").Value = 167994
$ws.Range("N84").Value = -177978
$ws.Range("H100").Value = 17500
$ws.Range("J100").Value = 17500
$ws.Range("L100").Value = 17500
$ws.Range("N100").Value = -19664
$ws.Range("H132").Value = 4172.933
$ws.Range("J132").Value = 10000
$ws.Range("L132").Value = 30000
$ws.Range("N132").Value = -35060

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 5080
$ws.Range("I29").Value = 5080
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 5080
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -4791
$ws.Range("N29").Value = $null
$ws.Range("H36").Value = 1000
$ws.Range("I36").Value = 1000
$ws.Range("K36").Value = 1000
$ws.Range("M36").Value = -466
$ws.Range("H86").Value = 13478.308
$ws.Range("I86").Value = 13072.708
$ws.Range("J86").Value = 14127.267
$ws.Range("K86").Value = 13072.708
$ws.Range("L86").Value = 14127.267
$ws.Range("M86").Value = -11949.708
$ws.Range("N86").Value = -16373.267
$ws.Range("H89").Value = 13478.308
$ws.Range("I89").Value = 13072.708
$ws.Range("J89").Value = 14127.267
$ws.Range("K89").Value = 65363.54
$ws.Range("L89").Value = 70636.33499999999
$ws.Range("M89").Value = -59747.54
$ws.Range("N89").Value = -81868.33499999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1377.25
$ws.Range("I31").Value = 1310.5454
$ws.Range("K31").Value = 1310.5454
$ws.Range("M31").Value = -1015.5454
$ws.Range("H34").Value = 1377.25
$ws.Range("I34").Value = 1310.5454
$ws.Range("K34").Value = 1310.5454
$ws.Range("M34").Value = -1108.5454
$ws.Range("H59").Value = 138333.33
$ws.Range("J59").Value = 146363.64
$ws.Range("L59").Value = 146363.64
$ws.Range("N59").Value = -148653.64
$ws.Range("H134").Value = 2306.1667
$ws.Range("I134").Value = 2524.5
$ws.Range("J134").Value = 1869.5
$ws.Range("K134").Value = 7573.5
$ws.Range("L134").Value = 5608.5
$ws.Range("M134").Value = -5038.5
$ws.Range("N134").Value = -10678.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4328.3335
$ws.Range("J80").Value = 4328.3335
$ws.Range("L80").Value = 12985.0005
$ws.Range("N80").Value = -14857.0005
$ws.Range("H83").Value = 4328.3335
$ws.Range("J83").Value = 4328.3335
$ws.Range("L83").Value = 38955.0015
$ws.Range("N83").Value = -48315.0015
$ws.Range("H93").Value = 19166.5
$ws.Range("J93").Value = 19166.5
$ws.Range("L93").Value = 57499.5
$ws.Range("N93").Value = -61243.5
$ws.Range("H122").Value = 844.5
$ws.Range("J122").Value = 824.75
$ws.Range("L122").Value = 7422.75
$ws.Range("N122").Value = -12322.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 1877.2
$ws.Range("I9").Value = 1877.2
$ws.Range("K9").Value = 1877.2
$ws.Range("M9").Value = -1707.2
$ws.Range("H44").Value = 19000
$ws.Range("I44").Value = 19000
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 19000
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -18404
$ws.Range("N44").Value = $null
$ws.Range("H80").Value = 2041.5
$ws.Range("I80").Value = 1449.8
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 1449.8
$ws.Range("L80").Value = 5000
$ws.Range("M80").Value = -451.8
$ws.Range("N80").Value = -6996
$ws.Range("H83").Value = 2041.5
$ws.Range("I83").Value = 1449.8
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 7249
$ws.Range("L83").Value = 25000
$ws.Range("M83").Value = -2257
$ws.Range("N83").Value = -34984
$ws.Range("H95").Value = 9999.666999999999
$ws.Range("J95").Value = 9999.5
$ws.Range("L95").Value = 9999.5
$ws.Range("N95").Value = -15491.5
$ws.Range("H132").Value = 15876952
$ws.Range("I132").Value = 4031.7144
$ws.Range("J132").Value = 47622790
$ws.Range("K132").Value = 12095.1432
$ws.Range("L132").Value = 142868370
$ws.Range("M132").Value = -9565.143199999999
$ws.Range("N132").Value = -142873430
$ws.Range("H136").Value = 68990.2
$ws.Range("J136").Value = 68990.2
$ws.Range("L136").Value = 206970.6
$ws.Range("N136").Value = -212070.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = $null
$ws.Range("H22").Value = 71430930
$ws.Range("I22").Value = 2000
$ws.Range("K22").Value = 2000
$ws.Range("M22").Value = -1705
$ws.Range("H25").Value = 5750
$ws.Range("J25").Value = 8000
$ws.Range("L25").Value = 8000
$ws.Range("N25").Value = -8460
$ws.Range("H26").Value = 14999.5
$ws.Range("J26").Value = 20000
$ws.Range("L26").Value = 20000
$ws.Range("N26").Value = -20590
$ws.Range("H27").Value = 71430930
$ws.Range("I27").Value = 2000
$ws.Range("K27").Value = 2000
$ws.Range("M27").Value = -1893
$ws.Range("H40").Value = 71431736
$ws.Range("I40").Value = 76925730
$ws.Range("K40").Value = 76925730
$ws.Range("M40").Value = -76925594
$ws.Range("H55").Value = 1619.9231
$ws.Range("I55").Value = 1419.75
$ws.Range("J55").Value = 1708.8889
$ws.Range("K55").Value = 1419.75
$ws.Range("L55").Value = 1708.8889
$ws.Range("M55").Value = -1246.75
$ws.Range("N55").Value = -2054.8889
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = $null
$ws.Range("H132").Value = 3444.75
$ws.Range("J132").Value = 3318.5625
$ws.Range("L132").Value = 9955.6875
$ws.Range("N132").Value = -15015.6875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 34750
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").Value = $null
$ws.Range("H107").Value = 4672
$ws.Range("I107").Value = 4672
$ws.Range("K107").Value = 14016
$ws.Range("M107").Value = -12096
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = $null
$ws.Range("H132").Value = 43488464
$ws.Range("I132").Value = 11455.2
$ws.Range("J132").Value = 333335200
$ws.Range("K132").Value = 34365.60000000001
$ws.Range("L132").Value = 1000005600
$ws.Range("M132").Value = -31835.60000000001
$ws.Range("N132").Value = -1000010660
